# Remove the four standalone paragraphs "Nf", "Nd", "Ndfz", "Nz" entirely
# (text + paragraph mark), and remove just the "nzffnfdd" run from the
# paragraph that follows them, leaving that paragraph (and its
# bookmarkStart/bookmarkEnd for "_GoBack") intact.

$d = $word.ActiveDocument

# Delete paragraphs from the end backward so earlier indices stay valid.
# Paragraph 4 = "Nf", 5 = "Nd", 6 = "Ndfz", 7 = "Nz".
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(6).Range.Delete()
$d.Paragraphs.Item(5).Range.Delete()
$d.Paragraphs.Item(4).Range.Delete()

# Remove just the "nzffnfdd" text, leaving the now-empty paragraph (with
# the _GoBack bookmark) in place.
[void]$d.Content.Find.Execute("nzffnfdd", $false, $false, $false, $false,
                               $false, $true, 1, $false, "", 2)
